$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 83
$ws.Range("E3").Value = 3.276141612960598
$ws.Range("E4").Value = 0.8515196826842663
$ws.Range("E6").Value = 2.376860703010635
$ws.Range("E7").Value = 3.306936630275666
